{"js": "const table = context.document.body.tables.getFirst();\ntable.load('values');\nawait context.sync();\n\nconst updates = [\n  { row: 0, col: 0, before: \"53\u00f75=10, 3\", after: \"25\u00f73=8, 1\" },\n  { row: 0, col: 1, before: \"27\u00f75=5, 2\", after: \"35\u00f79=3, 8\" },\n  { row: 0, col: 2, before: \"18\u00f77=2, 4\", after: \"55\u00f77=7, 6\" },\n  { row: 0, col: 3, before: \"43\u00f76=7, 1\", after: \"34\u00f74=8, 2\" },\n  { row: 0, col: 4, before: \"79\u00f78=9, 7\", after: \"41\u00f72=20, 1\" },\n  { row: 4, col: 0, before: \"83\u00f78=10, 3\", after: \"48\u00f74=12, 0\" },\n  { row: 4, col: 1, before: \"88\u00f74=22, 0\", after: \"18\u00f78=2, 2\" },\n  { row: 4, col: 2, before: \"46\u00f79=5, 1\", after: \"18\u00f78=2, 2\" },\n  { row: 4, col: 3, before: \"78\u00f77=11, 1\", after: \"28\u00f79=3, 1\" },\n  { row: 4, col: 4, before: \"55\u00f79=6, 1\", after: \"28\u00f79=3, 1\" },\n  { row: 8, col: 0, before: \"86\u00f75=17, 1\", after: \"89\u00f77=12, 5\" },\n  { row: 8, col: 1, before: \"89\u00f75=17, 4\", after: \"67\u00f76=11, 1\" },\n  { row: 8, col: 2, before: \"97\u00f73=32, 1\", after: \"19\u00f75=3, 4\" },\n  { row: 8, col: 3, before: \"84\u00f72=42, 0\", after: \"66\u00f73=22, 0\" },\n  { row: 8, col: 4, before: \"78\u00f77=11, 1\", after: \"89\u00f72=44, 1\" },\n  { row: 12, col: 0, before: \"41\u00f79=4, 5\", after: \"52\u00f73=17, 1\" },\n  { row: 12, col: 1, before: \"80\u00f77=11, 3\", after: \"52\u00f74=13, 0\" },\n  { row: 12, col: 2, before: \"96\u00f75=19, 1\", after: \"43\u00f74=10, 3\" },\n  { row: 12, col: 3, before: \"10\u00f79=1, 1\", after: \"65\u00f78=8, 1\" },\n  { row: 12, col: 4, before: \"55\u00f75=11, 0\", after: \"77\u00f76=12, 5\" },\n  { row: 16, col: 0, before: \"13\u00f77=1, 6\", after: \"54\u00f74=13, 2\" },\n  { row: 16, col: 1, before: \"33\u00f72=16, 1\", after: \"94\u00f73=31, 1\" },\n  { row: 16, col: 2, before: \"62\u00f78=7, 6\", after: \"14\u00f75=2, 4\" },\n  { row: 16, col: 3, before: \"83\u00f75=16, 3\", after: \"15\u00f77=2, 1\" },\n  { row: 16, col: 4, before: \"38\u00f74=9, 2\", after: \"60\u00f78=7, 4\" },\n];\n\nfor (const u of updates) {\n  const current = table.values[u.row][u.col];\n  if (current !== u.before) {\n    throw new Error(`Unexpected cell text at [${u.row},${u.col}]: ${JSON.stringify(current)} (expected ${JSON.stringify(u.before)})`);\n  }\n  table.getCell(u.row, u.col).value = u.after;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$t = $d.Tables(1)\n\n$updates = @(\n  @{row=1; col=1; before=\"53\u00f75=10, 3\"; after=\"25\u00f73=8, 1\"}\n  @{row=1; col=2; before=\"27\u00f75=5, 2\"; after=\"35\u00f79=3, 8\"}\n  @{row=1; col=3; before=\"18\u00f77=2, 4\"; after=\"55\u00f77=7, 6\"}\n  @{row=1; col=4; before=\"43\u00f76=7, 1\"; after=\"34\u00f74=8, 2\"}\n  @{row=1; col=5; before=\"79\u00f78=9, 7\"; after=\"41\u00f72=20, 1\"}\n  @{row=5; col=1; before=\"83\u00f78=10, 3\"; after=\"48\u00f74=12, 0\"}\n  @{row=5; col=2; before=\"88\u00f74=22, 0\"; after=\"18\u00f78=2, 2\"}\n  @{row=5; col=3; before=\"46\u00f79=5, 1\"; after=\"18\u00f78=2, 2\"}\n  @{row=5; col=4; before=\"78\u00f77=11, 1\"; after=\"28\u00f79=3, 1\"}\n  @{row=5; col=5; before=\"55\u00f79=6, 1\"; after=\"28\u00f79=3, 1\"}\n  @{row=9; col=1; before=\"86\u00f75=17, 1\"; after=\"89\u00f77=12, 5\"}\n  @{row=9; col=2; before=\"89\u00f75=17, 4\"; after=\"67\u00f76=11, 1\"}\n  @{row=9; col=3; before=\"97\u00f73=32, 1\"; after=\"19\u00f75=3, 4\"}\n  @{row=9; col=4; before=\"84\u00f72=42, 0\"; after=\"66\u00f73=22, 0\"}\n  @{row=9; col=5; before=\"78\u00f77=11, 1\"; after=\"89\u00f72=44, 1\"}\n  @{row=13; col=1; before=\"41\u00f79=4, 5\"; after=\"52\u00f73=17, 1\"}\n  @{row=13; col=2; before=\"80\u00f77=11, 3\"; after=\"52\u00f74=13, 0\"}\n  @{row=13; col=3; before=\"96\u00f75=19, 1\"; after=\"43\u00f74=10, 3\"}\n  @{row=13; col=4; before=\"10\u00f79=1, 1\"; after=\"65\u00f78=8, 1\"}\n  @{row=13; col=5; before=\"55\u00f75=11, 0\"; after=\"77\u00f76=12, 5\"}\n  @{row=17; col=1; before=\"13\u00f77=1, 6\"; after=\"54\u00f74=13, 2\"}\n  @{row=17; col=2; before=\"33\u00f72=16, 1\"; after=\"94\u00f73=31, 1\"}\n  @{row=17; col=3; before=\"62\u00f78=7, 6\"; after=\"14\u00f75=2, 4\"}\n  @{row=17; col=4; before=\"83\u00f75=16, 3\"; after=\"15\u00f77=2, 1\"}\n  @{row=17; col=5; before=\"38\u00f74=9, 2\"; after=\"60\u00f78=7, 4\"}\n)\n\nforeach ($u in $updates) {\n  $cell = $t.Cell($u.row, $u.col)\n  if (-not $cell.Range.Text.StartsWith($u.before)) {\n    throw \"Unexpected cell text at row=$($u.row) col=$($u.col): $($cell.Range.Text)\"\n  }\n  $cell.Range.Text = $u.after\n}\n"}
